# Version 000 - Commit and Push na Branch correta
#
# This script reproduces the edits made to resultado_organizacional.xlsx:
#  - Insert a new "Total_Membros" column (right after "Entropia Total")
#  - Insert two new columns "Org_Size" and "CC" (right before "Valor_da_Estrutura")
#  - Append two new columns "Complexidade Organizacional" and "Qualidade Prevista"
#  - Drop the per-cell/per-column custom formatting on the original data columns
#  - Apply a 3-decimal numeric format to "Valor_da_Estrutura"/"Complexidade Organizacional"
#  - Apply a percentage format to "Qualidade Prevista"
#  - Change the header row vertical alignment from centered to top
#  - Leave the sheet with columns B:D selected (as the editor last left it)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlTop    = -4160

# ---------------------------------------------------------------------------
# 1) Insert "Total_Membros" column right after "Entropia Total" (i.e. before
#    the existing "Niveis" column, which is column E before this edit).
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).Insert()
$ws.Range("E1").Value = "Total_Membros"
$ws.Range("E2").Value = 10
$ws.Range("E3").Value = 4

# ---------------------------------------------------------------------------
# 2) Insert "Org_Size" and "CC" columns right before "Valor_da_Estrutura"
#    (which is now column K after the previous insert).
# ---------------------------------------------------------------------------
$ws.Columns.Item(11).Insert()
$ws.Columns.Item(11).Insert()
$ws.Range("K1").Value = "Org_Size"
$ws.Range("L1").Value = "CC"
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 200
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 2000

# ---------------------------------------------------------------------------
# 3) Append "Complexidade Organizacional" and "Qualidade Prevista" columns
#    after "Valor_da_Estrutura" (now column M).
# ---------------------------------------------------------------------------
$ws.Range("N1").Value = "Complexidade Organizacional"
$ws.Range("O1").Value = "Qualidade Prevista"
$ws.Range("N2").Value = 103
$ws.Range("O2").Value = 0.96153846153846156
$ws.Range("N3").Value = 1001.2
$ws.Range("O3").Value = 0.099780482937537412

# ---------------------------------------------------------------------------
# 4) Strip the old centered / 2-decimal formatting from every data column
#    (A:L). They go back to the workbook's plain "Normal" style.
# ---------------------------------------------------------------------------
$ws.Range("A2:L3").Style = "Normal"

# ---------------------------------------------------------------------------
# 5) New numeric formatting for the value columns:
#    - M:N ("Valor_da_Estrutura", "Complexidade Organizacional") -> 0.000
#    - O   ("Qualidade Prevista") -> percentage, using the built-in
#      "Percent" cell style, centered.
# ---------------------------------------------------------------------------
$ws.Range("M2:N3").NumberFormat = "0.000"
$ws.Range("M2:N3").HorizontalAlignment = $xlCenter
$ws.Range("M2:N3").VerticalAlignment = $xlCenter

$ws.Range("O2:O3").Style = "Percent"
$ws.Range("O2:O3").NumberFormat = "0%"
$ws.Range("O2:O3").HorizontalAlignment = $xlCenter
$ws.Range("O2:O3").VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# 6) Header row: keep bold / border / centered horizontally, but switch the
#    vertical alignment from centered to top.
# ---------------------------------------------------------------------------
$ws.Range("A1:O1").VerticalAlignment = $xlTop

# ---------------------------------------------------------------------------
# 7) Leave the selection the way the author last left it.
# ---------------------------------------------------------------------------
$ws.Range("B1:D1048576").Select()
